# CLI_Data_Cleaning_Tool/output.xlsx - test common cleaning functions
# Clean up header/data casing and remove the duplicated sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the repeated "john doe" rows (4-13), keeping the header + first row.
$ws.Range("A4:C13").EntireRow.Delete()

# Normalize the header row to title case.
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Email"

# Normalize the remaining data row's name to title case (email stays as-is).
$ws.Range("B3").Value = "John Doe"
$ws.Range("C3").Value = "johndoe@gmail.com"

# Keep A1 registered as part of the sheet's used range (matches original
# layout where row 1 / column A were present but blank) so the sheet
# dimension collapses to A1:C3 instead of B2:C3.
$ws.Range("A1").Style = "Normal"
